$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

foreach ($col in $cols) {
    $cell13 = $ws.Range($col + "13")
    $cell14 = $ws.Range($col + "14")
    $v13 = $cell13.Value2
    $v14 = $cell14.Value2
    $cell13.Value2 = $v14
    $cell14.Value2 = $v13
}
